$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Existing layout:  A1=Rank  B1=Team  C1=ExpPoints
# New layout:       A1=Rank  B1=Team  C1=WIN  D1=TOP4  E1=TOP5  F1=RELEGATION  G1=ExpPoints
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# Give the new header cells the same look (bold / bordered / centered) as
# the existing header cells by copying the formatting from B1, which
# already carries that style - this reuses the existing style record
# instead of fabricating a new one.
$ws.Range("B1").Copy()
$ws.Range("D1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Team / ExpPoints data --------------------------------------------
# New simulated ExpPoints values (column G), keyed by team name.
$expPoints = @{
    "Arsenal"                 = 77.78415282723644
    "Manchester City"         = 72.09722680216554
    "Liverpool"               = 68.60549704904238
    "Chelsea"                 = 61.52722018553493
    "Aston Villa"             = 58.88880445799801
    "Crystal Palace"          = 57.98105521406693
    "Newcastle United"        = 57.42532546530259
    "Brighton & Hove Albion"  = 54.69731458139393
    "Tottenham Hotspur"       = 52.14388024299829
    "AFC Bournemouth"         = 52.06710236435932
    "Manchester United"       = 51.28649867493214
    "Brentford"               = 51.10735898605598
    "Everton"                 = 45.64404316838646
    "Fulham"                  = 43.55279677261844
    "Nottingham Forest"       = 40.70418084528933
    "West Ham United"         = 38.37735420158572
    "Sunderland"              = 37.63340043412545
    "Burnley"                 = 35.14843568014878
    "Leeds United"            = 34.80902214506362
    "Wolverhampton Wanderers" = 33.00560166274826
}

# New team order for rows 2..21 (ranks in column A stay 1..20 in row order;
# Aston Villa / Newcastle United and Tottenham Hotspur / AFC Bournemouth
# swap places relative to the previous table).
$teamOrder = @(
    "Arsenal",
    "Manchester City",
    "Liverpool",
    "Chelsea",
    "Aston Villa",
    "Crystal Palace",
    "Newcastle United",
    "Brighton & Hove Albion",
    "Tottenham Hotspur",
    "AFC Bournemouth",
    "Manchester United",
    "Brentford",
    "Everton",
    "Fulham",
    "Nottingham Forest",
    "West Ham United",
    "Sunderland",
    "Burnley",
    "Leeds United",
    "Wolverhampton Wanderers"
)

for ($i = 0; $i -lt $teamOrder.Length; $i++) {
    $row = $i + 2
    $team = $teamOrder[$i]

    $ws.Cells.Item($row, 2).Value = $team

    # WIN / TOP4 / TOP5 / RELEGATION (columns C-F) are placeholders for the
    # upcoming Monte Carlo simulation - leave them blank for now, clearing
    # out whatever used to live in column C (the old ExpPoints number).
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 3).Style = "Normal"
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 4).Style = "Normal"
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 5).Style = "Normal"
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 6).Style = "Normal"

    # ExpPoints now lives in column G with the refreshed projection.
    $ws.Cells.Item($row, 7).Value = $expPoints[$team]
}
